# Add a new "Assignment_16" row (row 17) to the Assignments tracker,
# matching the layout/format of the existing rows (e.g. row 16), with
# its GitHub link hyperlinked, and move the selection to E17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the assignment name first so it lands at the next shared-string
# slot before the hyperlink URL text does.
$ws.Range("A17").Value = "Assignment_16"

# Create the hyperlink on B17 (this also writes the display text).
$ws.Hyperlinks.Add($ws.Range("B17"), "https://github.com/Vasanth30e/Assignments_Phase2/tree/master/Assignment_16/TaskCRUD")

# Copy the row-16 formatting (fonts/borders/number format/row height)
# down into row 17 so the new row looks like the others.
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Rows(17).RowHeight = 34.5

# Re-assert the cell values/date (PasteSpecial only touched formatting).
$ws.Range("A17").Value = "Assignment_16"
$ws.Range("B17").Value = "https://github.com/Vasanth30e/Assignments_Phase2/tree/master/Assignment_16/TaskCRUD"
$ws.Range("C17").Value = 45170

# Match the saved selection from the edit (cell E17).
[void]$ws.Range("E17").Select()
